# The edit removes the original row 76 (SOLEVUL / -1 / -1 / I, with the
# "NA w inferred l inferred" note) from the Bio_Trust_24 sheet. All the
# subsequent rows (formerly 77-97) shift up by one row, which is exactly
# what deleting the worksheet row does, and the sheet's used range/
# dimension shrinks from A1:Q97 to A1:Q96 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(76).Delete()
